$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Columns.Item(3).Delete()
$ws1.Range("C2:C4").ClearContents()
$ws1.Range("C2").Select()

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("A1").Value = "0,0"
$ws2.Range("B1").Value = "0,1"
$ws2.Range("C1").Value = "0,2"
$ws2.Range("D1").Value = "0,3"
$ws2.Range("A1:D1").Font.Bold = $true

$ws2.Range("A2").Value = "1,1"
$ws2.Range("B2").Value = "1,2"
$ws2.Range("C2").Value = "1,3"
$ws2.Range("D2").Value = "1,4"

$ws2.Range("A3").Value = "locked_out_user"
$ws2.Range("B3").Value = "secret_sauce"
$ws2.Range("D3").Value = "Fail"

$ws2.Range("A4").Value = "problem_user"
$ws2.Range("B4").Value = "secret_sauce"
$ws2.Range("D4").Value = "Pass"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 14.0
$ws2.Columns.Item(2).ColumnWidth = 11.2
$ws2.Columns.Item(3).ColumnWidth = 8.0
$ws2.Columns.Item(4).ColumnWidth = 6.0

$ws2.Range("B2").Select()
$excel.ActiveWindow.Zoom = 120

# Re-activate the Credentials sheet so it stays the selected tab.
$ws1.Activate()
$ws1.Range("C2").Select()

Write-Host "done"
